$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data value in C5 (2000 -> 1000)
$ws.Range("C5").Value = 1000

# Update the active selection to C5 (was C6)
$ws.Range("C5").Select()
